# Workbook/sheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 14-20 newly assigned to "Trí" ---
# (assigned first so the rebuilt shared-string table keeps "Trí" at the
#  lower index, matching the original "Hieu" slot being renamed)
$ws.Range("E14").Value2 = "Trí"
$ws.Range("E15").Value2 = "Trí"
$ws.Range("E16").Value2 = "Trí"
$ws.Range("E17").Value2 = "Trí"
$ws.Range("E18").Value2 = "Trí"
$ws.Range("E19").Value2 = "Trí"
$ws.Range("E20").Value2 = "Trí"

# --- "first part" rows re-assigned from "Hieu" to "Hiếu" ---
$ws.Range("E2").Value2  = "Hiếu"
$ws.Range("D8").Value2  = 1
$ws.Range("E8").Value2  = "Hiếu"
$ws.Range("E9").Value2  = "Hiếu"
$ws.Range("E10").Value2 = "Hiếu"
$ws.Range("E11").Value2 = "Hiếu"

# --- Rows 12 & 13 newly filled in (dates + 100% + "Hiếu") ---
$ws.Range("B12").Value2 = 45693
$ws.Range("C12").Value2 = 45693
$ws.Range("D12").Value2 = 1
$ws.Range("E12").Value2 = "Hiếu"

$ws.Range("B13").Value2 = 45693
$ws.Range("C13").Value2 = 45693
$ws.Range("D13").Value2 = 1
$ws.Range("E13").Value2 = "Hiếu"

# --- selection moves to F9 ---
$ws.Range("F9").Select() | Out-Null
